# "Min model size ran" - fill in the "Min Input Space NN" (row 19) and
# "Min Input Space NN Regr." (row 22) result rows, center the data-grid
# cells, set the page to portrait, and move the selection to H22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New results for row 19 (Min Input Space NN: class accuracy) ---
$ws.Range("B19").Value = 0.92025518417358398
$ws.Range("C19").Value = 0.93033401171366303
$ws.Range("D19").Value = 0.92022582888603199
$ws.Range("E19").Value = 0.96416300535202004
$ws.Range("F19").Value = 0.95139908790588301
$ws.Range("G19").Value = 0.97539975245793598
$ws.Range("H19").Value = "63 params"
$ws.Range("K19").Value = 0.88039215405782001
$ws.Range("L19").Value = 0.89313725630442298
$ws.Range("M19").Value = 0.86764705181121804
$ws.Range("N19").Value = 0.88725491364796905
$ws.Range("O19").Value = 0.87352943420410101
$ws.Range("P19").Value = 0.83627450466155995

# --- New results for row 22 (Min Input Space NN Regr.: mean squared error) ---
$ws.Range("B22").Value = 0.812457581361134
$ws.Range("C22").Value = 0.72090607881545998
$ws.Range("D22").Value = 0.80801308155059803
$ws.Range("E22").Value = 0.619081750512123
$ws.Range("F22").Value = 0.32622307538986201
$ws.Range("G22").Value = 0.29697552323341297
$ws.Range("K22").Value = 2.8074042797088601
$ws.Range("L22").Value = 3.0053186416625901
$ws.Range("M22").Value = 2.9003276824951101
$ws.Range("N22").Value = 3.4158694744110099
$ws.Range("O22").Value = 3.3557066917419398
$ws.Range("P22").Value = 4.5839817523956299

# --- Center-align the whole results grid (rows 9-25) in both tables ---
$ws.Range("B9:F25").HorizontalAlignment = -4108
$ws.Range("K9:P25").HorizontalAlignment = -4108
$ws.Range("G9:G25").HorizontalAlignment = -4108

# --- Page orientation: portrait ---
$ws.PageSetup.Orientation = 1

# --- Move the active selection to H22 ---
[void]$ws.Range("H22").Select()
